$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C72").Value = "S220"
$ws.Range("A72").Value = "Look Into the Future"
$ws.Range("B72").Value = "Generic"
$ws.Range("D72").Value = "Magic"
$ws.Range("E72").Value = "Rare"

$ws.Range("C73").Value = "S221"
$ws.Range("A73").Value = "Ice Mammoth"
$ws.Range("B73").Value = "Generic"
$ws.Range("D73").Value = "Creature"
$ws.Range("E73").Value = "Rare"

$ws.Range("C74").Value = "S222"
$ws.Range("A74").Value = "Necromancy"
$ws.Range("B74").Value = "Generic"
$ws.Range("D74").Value = "Magic"
$ws.Range("E74").Value = "Rare"

[void]$ws.Range("A75").Select()
